$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 146.23077
$ws.Range("I9").Value = 151.75
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 151.75
$ws.Range("L9").Value = 80
$ws.Range("M9").Value = 17.25
$ws.Range("N9").Value = -418
$ws.Range("H21").Value = 11636.272
$ws.Range("I21").Value = 8000
$ws.Range("J21").Value = 15999.8
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 15999.8
$ws.Range("M21").Value = -7532
$ws.Range("H23").Value = 11636.272
$ws.Range("I23").Value = 8000
$ws.Range("J23").Value = 15999.8
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 15999.8
$ws.Range("M23").Value = -7766
$ws.Range("H38").Value = 3305.2222
$ws.Range("I38").Value = 172
$ws.Range("J38").Value = 5299.091
$ws.Range("K38").Value = 516
$ws.Range("L38").Value = 15897.273
$ws.Range("M38").Value = -144
$ws.Range("H62").Value = 1293.5454
$ws.Range("I62").Value = 1516.125
$ws.Range("J62").Value = 700
$ws.Range("K62").Value = 1516.125
$ws.Range("L62").Value = 700
$ws.Range("M62").Value = -892.125
$ws.Range("N62").Value = -1948
$ws.Range("H65").Value = 1293.5454
$ws.Range("I65").Value = 1516.125
$ws.Range("J65").Value = 700
$ws.Range("K65").Value = 7580.625
$ws.Range("L65").Value = 3500
$ws.Range("M65").Value = -4460.625
$ws.Range("N65").Value = -9740
$ws.Range("H138").Value = 2810.9365
$ws.Range("I138").Value = 1191.75
$ws.Range("J138").Value = 3046.4546
$ws.Range("K138").Value = 3575.25
$ws.Range("L138").Value = 9139.363799999999
$ws.Range("M138").Value = 1564.75
$ws.Range("N138").Value = -19419.3638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 737
$ws.Range("I2").Value = 605.5
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 605.5
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -492.5
$ws.Range("N2").Value = -1226
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H23").Value = 25083.666
$ws.Range("I23").Value = 25000
$ws.Range("J23").Value = 25251
$ws.Range("K23").Value = 25000
$ws.Range("L23").Value = 25251
$ws.Range("M23").Value = -24741
$ws.Range("N23").Value = -25769
$ws.Range("H37").Value = 25185
$ws.Range("I37").Value = 12100
$ws.Range("J37").Value = 30091.875
$ws.Range("K37").Value = 12100
$ws.Range("L37").Value = 30091.875
$ws.Range("M37").Value = -11827
$ws.Range("N37").Value = -30637.875
$ws.Range("H44").Value = 35651.547
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 35651.547
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 35651.547
$ws.Range("N44").Value = -36627.547
$ws.Range("H55").Value = 25174.445
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 25174.445
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 25174.445
$ws.Range("N55").Value = -25804.445
$ws.Range("H80").Value = 37807.11
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 37807.11
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 37807.11
$ws.Range("N80").Value = -39803.11
$ws.Range("H83").Value = 37807.11
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 37807.11
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 113421.33
$ws.Range("N83").Value = -123405.33
$ws.Range("H116").Value = 737
$ws.Range("I116").Value = 605.5
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 605.5
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1688.5
$ws.Range("N116").Value = -5588
$ws.Range("H132").Value = 3454.7646
$ws.Range("I132").Value = 1900.7142
$ws.Range("J132").Value = 4542.6
$ws.Range("K132").Value = 5702.142599999999
$ws.Range("L132").Value = 13627.8
$ws.Range("M132").Value = -3172.142599999999
$ws.Range("N132").Value = -18687.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 737
$ws.Range("I3").Value = 605.5
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 605.5
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -491.5
$ws.Range("N3").Value = -1228
$ws.Range("H22").Value = 176.85715
$ws.Range("I22").Value = 132.33333
$ws.Range("J22").Value = 444
$ws.Range("K22").Value = 132.33333
$ws.Range("L22").Value = 444
$ws.Range("M22").Value = 40.66667000000001
$ws.Range("N22").Value = -790

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5649.1113
$ws.Range("I31").Value = 2313.4688
$ws.Range("J31").Value = 10500.954
$ws.Range("K31").Value = 2313.4688
$ws.Range("L31").Value = 10500.954
$ws.Range("M31").Value = -2018.4688
$ws.Range("H34").Value = 5649.1113
$ws.Range("I34").Value = 2313.4688
$ws.Range("J34").Value = 10500.954
$ws.Range("K34").Value = 2313.4688
$ws.Range("L34").Value = 10500.954
$ws.Range("M34").Value = -2111.4688
$ws.Range("H58").Value = 2272.8298
$ws.Range("I58").Value = 1411.75
$ws.Range("J58").Value = 5090.909
$ws.Range("K58").Value = 1411.75
$ws.Range("L58").Value = 5090.909
$ws.Range("M58").Value = -1208.75
$ws.Range("H98").Value = 41161.5
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 41161.5
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 41161.5
$ws.Range("N98").Value = -45653.5
$ws.Range("H136").Value = 2272.8298
$ws.Range("I136").Value = 1411.75
$ws.Range("J136").Value = 5090.909
$ws.Range("K136").Value = 4235.25
$ws.Range("L136").Value = 15272.727
$ws.Range("M136").Value = -1685.25
$ws.Range("H137").Value = 31230
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 31230
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 31230
$ws.Range("N137").Value = -41430

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 14746.444
$ws.Range("I34").Value = 16967.143
$ws.Range("J34").Value = 13333.272
$ws.Range("K34").Value = 50901.429
$ws.Range("L34").Value = 39999.81600000001
$ws.Range("M34").Value = -50817.429
$ws.Range("N34").Value = -40167.81600000001
$ws.Range("H39").Value = 13399.55
$ws.Range("I39").Value = 4999
$ws.Range("J39").Value = 13841.685
$ws.Range("K39").Value = 14997
$ws.Range("L39").Value = 41525.055
$ws.Range("M39").Value = -14703
$ws.Range("N39").Value = -42113.055
$ws.Range("H55").Value = 4226.25
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 5301.6665
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 15904.9995
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -16258.9995
$ws.Range("H68").Value = 621.25
$ws.Range("I68").Value = 596
$ws.Range("J68").Value = 663.3333
$ws.Range("K68").Value = 1788
$ws.Range("L68").Value = 1989.9999
$ws.Range("M68").Value = -977
$ws.Range("N68").Value = -3611.9999
$ws.Range("H71").Value = 621.25
$ws.Range("I71").Value = 596
$ws.Range("J71").Value = 663.3333
$ws.Range("K71").Value = 5364
$ws.Range("L71").Value = 5969.9997
$ws.Range("M71").Value = -1308
$ws.Range("N71").Value = -14081.9997
$ws.Range("H107").Value = 23857676
$ws.Range("I107").Value = 435.875
$ws.Range("J107").Value = 38539056
$ws.Range("K107").Value = 1307.625
$ws.Range("L107").Value = 115617168
$ws.Range("M107").Value = 612.375
$ws.Range("N107").Value = -115621008
$ws.Range("H113").Value = 5435357
$ws.Range("I113").Value = 634.1
$ws.Range("J113").Value = 9615913
$ws.Range("K113").Value = 1902.3
$ws.Range("L113").Value = 28847739
$ws.Range("M113").Value = 267.6999999999998
$ws.Range("N113").Value = -28852079
$ws.Range("H129").Value = 3603.3333
$ws.Range("I129").Value = 2803.75
$ws.Range("J129").Value = 10000
$ws.Range("K129").Value = 8411.25
$ws.Range("L129").Value = 30000
$ws.Range("M129").Value = -3411.25
$ws.Range("N129").Value = -40000
$ws.Range("H131").Value = 788.09
$ws.Range("I131").Value = 311.1111
$ws.Range("J131").Value = 835.26373
$ws.Range("K131").Value = 933.3333
$ws.Range("L131").Value = 2505.79119
$ws.Range("M131").Value = 4106.6667
$ws.Range("N131").Value = -12585.79119
$ws.Range("H132").Value = 2162.1614
$ws.Range("I132").Value = 997.9167
$ws.Range("J132").Value = 2897.4736
$ws.Range("K132").Value = 8981.2503
$ws.Range("L132").Value = 26077.2624
$ws.Range("M132").Value = -6451.2503
$ws.Range("N132").Value = -31137.2624

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3469.7
$ws.Range("I132").Value = 1876.6923
$ws.Range("J132").Value = 6428.143
$ws.Range("K132").Value = 5630.0769
$ws.Range("L132").Value = 19284.429
$ws.Range("M132").Value = -3100.0769
$ws.Range("N132").Value = -24344.429
$ws.Range("H135").Value = 32348.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 32348.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 32348.75
$ws.Range("N135").Value = -42488.75
$ws.Range("H137").Value = 43769.332
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 43769.332
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 43769.332
$ws.Range("N137").Value = -53969.332
$ws.Range("H138").Value = 42220
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 42220
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 42220
$ws.Range("N138").Value = -52500
$ws.Range("H140").Value = 39736.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 39736.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 39736.668
$ws.Range("N140").Value = -50096.668

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5741.4165
$ws.Range("I7").Value = 2074.75
$ws.Range("J7").Value = 7574.75
$ws.Range("K7").Value = 2074.75
$ws.Range("L7").Value = 7574.75
$ws.Range("M7").Value = -1962.75
$ws.Range("N7").Value = -7798.75
$ws.Range("H126").Value = 5741.4165
$ws.Range("I126").Value = 2074.75
$ws.Range("J126").Value = 7574.75
$ws.Range("K126").Value = 6224.25
$ws.Range("L126").Value = 22724.25
$ws.Range("M126").Value = -3754.25
$ws.Range("N126").Value = -27664.25
$ws.Range("H132").Value = 4373.231
$ws.Range("I132").Value = 3405.3684
$ws.Range("J132").Value = 7000.2856
$ws.Range("K132").Value = 10216.1052
$ws.Range("L132").Value = 21000.8568
$ws.Range("M132").Value = -7686.1052

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 15153794
$ws.Range("I132").Value = 1228.3
$ws.Range("J132").Value = 27780932
$ws.Range("K132").Value = 3684.9
$ws.Range("L132").Value = 83342796
$ws.Range("M132").Value = -1154.9
$ws.Range("N132").Value = -83347856
